$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated probability matrix values from games pulled march 7

# Row 2
$ws.Range("B2").Value = 0.1940298507462687
$ws.Range("C2").Value = 0.5671641791044776
$ws.Range("J2").Value = 0.01791044776119403
$ws.Range("P2").Value = 0.1402985074626866
$ws.Range("S2").Value = 0.08059701492537313

# Row 3
$ws.Range("B3").Value = 0.01515151515151515
$ws.Range("C3").Value = 0.02525252525252525
$ws.Range("J3").Value = 0.0303030303030303
$ws.Range("P3").Value = 0.7727272727272727
$ws.Range("S3").Value = 0.1565656565656566

# Row 6
$ws.Range("B6").Value = 0.08333333333333333
$ws.Range("D6").Value = 0.01666666666666667
$ws.Range("F6").Value = 0.075
$ws.Range("J6").Value = 0.2916666666666667
$ws.Range("O6").Value = 0.0125
$ws.Range("Q6").Value = 0.1916666666666667
$ws.Range("R6").Value = 0.05833333333333333
$ws.Range("S6").Value = 0.2708333333333333

# Row 7
$ws.Range("B7").Value = 0.137809187279152
$ws.Range("D7").Value = 0.01060070671378092
$ws.Range("F7").Value = 0.02826855123674912
$ws.Range("J7").Value = 0.1448763250883392
$ws.Range("O7").Value = 0.0176678445229682
$ws.Range("Q7").Value = 0.2084805653710247
$ws.Range("R7").Value = 0.09540636042402827
$ws.Range("S7").Value = 0.3568904593639576

# Row 8
$ws.Range("B8").Value = 0.07692307692307693
$ws.Range("D8").Value = 0.02439024390243903
$ws.Range("F8").Value = 0.04127579737335835
$ws.Range("J8").Value = 0.1144465290806754
$ws.Range("O8").Value = 0.01125703564727955
$ws.Range("Q8").Value = 0.2382739212007505
$ws.Range("R8").Value = 0.09380863039399624
$ws.Range("S8").Value = 0.399624765478424

# Row 9
$ws.Range("B9").Value = 0.09523809523809523
$ws.Range("D9").Value = 0.01984126984126984
$ws.Range("F9").Value = 0.07142857142857142
$ws.Range("J9").Value = 0.07539682539682539
$ws.Range("O9").Value = 0.0119047619047619
$ws.Range("Q9").Value = 0.1904761904761905
$ws.Range("R9").Value = 0.07539682539682539
$ws.Range("S9").Value = 0.4603174603174603

# Row 10
$ws.Range("B10").Value = 0.1004953998584572
$ws.Range("D10").Value = 0.01556970983722576
$ws.Range("F10").Value = 0.07006369426751592
$ws.Range("J10").Value = 0.1280962491153574
$ws.Range("O10").Value = 0.02193913658881812
$ws.Range("Q10").Value = 0.2469922151450814
$ws.Range("R10").Value = 0.08917197452229299
$ws.Range("S10").Value = 0.3276716206652512

# Row 11
$ws.Range("G11").Value = 0.1590361445783132
$ws.Range("J11").Value = 0.07228915662650602
$ws.Range("K11").Value = 0.1927710843373494
$ws.Range("L11").Value = 0.5614457831325301
$ws.Range("S11").Value = 0.01445783132530121

# Row 12
$ws.Range("F12").Value = 0.004098360655737705
$ws.Range("G12").Value = 0.7418032786885246
$ws.Range("J12").Value = 0.1680327868852459
$ws.Range("K12").Value = 0.02049180327868852
$ws.Range("L12").Value = 0.04098360655737705
$ws.Range("S12").Value = 0.02459016393442623

# Row 15
$ws.Range("F15").Value = 0.02531645569620253
$ws.Range("H15").Value = 0.1392405063291139
$ws.Range("I15").Value = 0.04641350210970464
$ws.Range("J15").Value = 0.3544303797468354
$ws.Range("K15").Value = 0.08438818565400844
$ws.Range("M15").Value = 0.01687763713080169
$ws.Range("O15").Value = 0.06329113924050633
$ws.Range("S15").Value = 0.270042194092827

# Row 16
$ws.Range("F16").Value = 0.01746724890829694
$ws.Range("H16").Value = 0.1965065502183406
$ws.Range("I16").Value = 0.1004366812227074
$ws.Range("J16").Value = 0.3799126637554585
$ws.Range("K16").Value = 0.1222707423580786
$ws.Range("M16").Value = 0.02620087336244541
$ws.Range("N16").Value = 0.008733624454148471
$ws.Range("O16").Value = 0.03056768558951965
$ws.Range("S16").Value = 0.1179039301310044

# Row 17
$ws.Range("F17").Value = 0.01762820512820513
$ws.Range("H17").Value = 0.1907051282051282
$ws.Range("I17").Value = 0.09935897435897435
$ws.Range("J17").Value = 0.4006410256410257
$ws.Range("K17").Value = 0.1105769230769231
$ws.Range("M17").Value = 0.02083333333333333
$ws.Range("N17").Value = 0.001602564102564103
$ws.Range("O17").Value = 0.0641025641025641
$ws.Range("S17").Value = 0.09455128205128205

# Row 18
$ws.Range("F18").Value = 0.01271186440677966
$ws.Range("H18").Value = 0.1991525423728814
$ws.Range("I18").Value = 0.08898305084745763
$ws.Range("J18").Value = 0.4110169491525424
$ws.Range("K18").Value = 0.0635593220338983
$ws.Range("M18").Value = 0.0211864406779661
$ws.Range("N18").Value = 0.00423728813559322
$ws.Range("O18").Value = 0.05932203389830509
$ws.Range("S18").Value = 0.1398305084745763

# Row 19
$ws.Range("F19").Value = 0.01460920379839299
$ws.Range("H19").Value = 0.2081811541271001
$ws.Range("I19").Value = 0.09861212563915267
$ws.Range("J19").Value = 0.3403944485025566
$ws.Range("K19").Value = 0.1344046749452155
$ws.Range("M19").Value = 0.02629656683710738
$ws.Range("N19").Value = 0.0007304601899196494
$ws.Range("O19").Value = 0.0606281957633309
$ws.Range("S19").Value = 0.1161431701972243
